# Working download and chatting functionality
# Rename sheets and replace their data tables.

$wb = $excel.ActiveWorkbook

# --- Rename sheets -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Football Goals"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Cricket Runs"

# --- Sheet 1: "Football Goals" -----------------------------------------
# New header row
$ws1.Cells.Item(1,1).Value = "Rank"
$ws1.Cells.Item(1,2).Value = "Player"
$ws1.Cells.Item(1,3).Value = "Team"
$ws1.Cells.Item(1,4).Value = "Season"
$ws1.Cells.Item(1,5).Value = "Goals"

# Copy the existing bold header style (from A1) onto the newly used C1:E1 cells
$ws1.Range("A1").Copy()
$ws1.Range("C1:E1").PasteSpecial(-4122)

$football = @(
    @(1, "Lionel Messi", "Barcelona", "2011/12", 73),
    @(2, "Ferenc Deak", "Szentlorinci", "1945/46", 66),
    @(2, "Gerd Muller", "Bayern Munich", "1972/73", 66),
    @(4, "Dixie Dean", "Everton", "1927/28", 63),
    @(5, "Cristiano Ronaldo", "Real Madrid", "2014/15", 61),
    @(6, "Cristiano Ronaldo", "Real Madrid", "2011/12", 60),
    @(6, "Lionel Messi", "Barcelona", "2012/13", 60),
    @(8, "Ferenc Deak", "Ferencvaros", "1948/49", 59),
    @(8, "Luis Suarez", "Barcelona", "2015/16", 59),
    @(10, "Lionel Messi", "Barcelona", "2014/15", 58)
)

$r = 2
foreach ($row in $football) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# --- Sheet 2: "Cricket Runs" --------------------------------------------
# New header row
$ws2.Cells.Item(1,1).Value = "Rank"
$ws2.Cells.Item(1,2).Value = "Player"
$ws2.Cells.Item(1,3).Value = "Country"
$ws2.Cells.Item(1,4).Value = "Runs"
$ws2.Cells.Item(1,5).Value = "Year"

# Copy the existing bold header style (from A1) onto the newly used D1:E1 cells
$ws2.Range("A1").Copy()
$ws2.Range("D1:E1").PasteSpecial(-4122)

$cricket = @(
    @(1, "KC Sangakkara", "SL", 2868, 2013),
    @(2, "RT Ponting", "AUS/ICC", 2833, 2005),
    @(3, "V Kohli", "IND", 2818, 2017),
    @(4, "V Kohli", "IND", 2735, 2018),
    @(5, "KS Williamson", "NZ", 2692, 2015)
)

$r = 2
foreach ($row in $cricket) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

Write-Host "Edits applied."
